$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16:XFD17").Select()
$ws.Range("A16:G17").EntireRow.Delete()
